# Update database table HoaDon
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Change C3 from 70 to 2
$ws.Range("C3").Value = 2

# Update product code SP0004 -> SP0006 (shared string used across column A)
$ws.Range("A3:A7").Value = "SP0006"

# Delete the last data row (row 7), which duplicated row 6
$ws.Rows.Item(7).Delete()

# Update selection to match target state
$ws.Range("D11").Select()
